{"js": "// Update the title date line and every arithmetic-problem cell in the\n// (single) table with the new values from the target revision.\n// Cell text is replaced row-by-row, left-to-right, matching document order.\n\nconst body = context.document.body;\n\n// --- 1) Title paragraph: \"2025-10-20 Monday\" -> \"2025-10-21 Tuesday\" ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length > 0) {\n  paragraphs.items[0].insertText(\"2025-10-21 Tuesday\", \"Replace\");\n}\n\n// --- 2) Table cells: replace every cell's text with the new value ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst newValues = [\n  [\"16+29=\", \"76-28=\", \"36+26=\", \"26+25=\", \"24+19=\"],\n  [\"36+58=\", \"71-56=\", \"48+33=\", \"19+6=\", \"9+34=\"],\n  [\"27+19=\", \"39+27=\", \"6+67=\", \"46+25=\", \"38+45=\"],\n  [\"5+39=\", \"20-12=\", \"57-18=\", \"40-13=\", \"27+54=\"],\n  [\"61-42=\", \"29+52=\", \"58+4=\", \"84-39=\", \"17+69=\"],\n  [\"63-6=\", \"39+24=\", \"15+36=\", \"94-35=\", \"91-83=\"],\n  [\"36+39=\", \"17+28=\", \"82-68=\", \"94-36=\", \"83-35=\"],\n  [\"69+7=\", \"22+49=\", \"27+7=\", \"91-53=\", \"72-53=\"],\n  [\"7+16=\", \"56-27=\", \"16+48=\", \"35+9=\", \"60-38=\"],\n  [\"84+9=\", \"37+18=\", \"26+29=\", \"65+9=\", \"42-15=\"],\n  [\"42-39=\", \"17+6=\", \"26+66=\", \"44-5=\", \"75-38=\"],\n  [\"18+54=\", \"13-9=\", \"33+8=\", \"41-12=\", \"35+7=\"],\n  [\"60-37=\", \"97-38=\", \"33+48=\", \"36+17=\", \"46-7=\"],\n  [\"47+37=\", \"76+18=\", \"28+47=\", \"37+34=\", \"67+15=\"],\n  [\"18+74=\", \"4+59=\", \"30-19=\", \"78-9=\", \"95-48=\"],\n  [\"52-29=\", \"18+66=\", \"8+15=\", \"67-29=\", \"50-44=\"],\n  [\"91-59=\", \"44+29=\", \"64-39=\", \"58-9=\", \"43-19=\"],\n  [\"91-74=\", \"45+48=\", \"54-46=\", \"73-9=\", \"96-88=\"],\n  [\"51-49=\", \"91-23=\", \"75-37=\", \"60-29=\", \"5+86=\"],\n  [\"86-37=\", \"94-25=\", \"80-73=\", \"68+7=\", \"58+35=\"],\n];\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n  table.values = newValues;\n}\n\nawait context.sync();\n", "ps1": "# Update the title date line and every arithmetic-problem cell in the\n# (single) table with the new values from the target revision.\n# Cell text is replaced row-by-row, left-to-right, matching document order.\n\n$d = $word.ActiveDocument\n\n# --- 1) Title paragraph: \"2025-10-20 Monday\" -> \"2025-10-21 Tuesday\" ---\n$d.Paragraphs.Item(1).Range.Text = \"2025-10-21 Tuesday\"\n\n# --- 2) Table cells: replace every cell's text with the new value ---\n$values = @(\n    @(\"16+29=\", \"76-28=\", \"36+26=\", \"26+25=\", \"24+19=\"),\n    @(\"36+58=\", \"71-56=\", \"48+33=\", \"19+6=\", \"9+34=\"),\n    @(\"27+19=\", \"39+27=\", \"6+67=\", \"46+25=\", \"38+45=\"),\n    @(\"5+39=\", \"20-12=\", \"57-18=\", \"40-13=\", \"27+54=\"),\n    @(\"61-42=\", \"29+52=\", \"58+4=\", \"84-39=\", \"17+69=\"),\n    @(\"63-6=\", \"39+24=\", \"15+36=\", \"94-35=\", \"91-83=\"),\n    @(\"36+39=\", \"17+28=\", \"82-68=\", \"94-36=\", \"83-35=\"),\n    @(\"69+7=\", \"22+49=\", \"27+7=\", \"91-53=\", \"72-53=\"),\n    @(\"7+16=\", \"56-27=\", \"16+48=\", \"35+9=\", \"60-38=\"),\n    @(\"84+9=\", \"37+18=\", \"26+29=\", \"65+9=\", \"42-15=\"),\n    @(\"42-39=\", \"17+6=\", \"26+66=\", \"44-5=\", \"75-38=\"),\n    @(\"18+54=\", \"13-9=\", \"33+8=\", \"41-12=\", \"35+7=\"),\n    @(\"60-37=\", \"97-38=\", \"33+48=\", \"36+17=\", \"46-7=\"),\n    @(\"47+37=\", \"76+18=\", \"28+47=\", \"37+34=\", \"67+15=\"),\n    @(\"18+74=\", \"4+59=\", \"30-19=\", \"78-9=\", \"95-48=\"),\n    @(\"52-29=\", \"18+66=\", \"8+15=\", \"67-29=\", \"50-44=\"),\n    @(\"91-59=\", \"44+29=\", \"64-39=\", \"58-9=\", \"43-19=\"),\n    @(\"91-74=\", \"45+48=\", \"54-46=\", \"73-9=\", \"96-88=\"),\n    @(\"51-49=\", \"91-23=\", \"75-37=\", \"60-29=\", \"5+86=\"),\n    @(\"86-37=\", \"94-25=\", \"80-73=\", \"68+7=\", \"58+35=\")\n)\n\n$t = $d.Tables.Item(1)\nfor ($r = 0; $r -lt $values.Count; $r++) {\n    for ($c = 0; $c -lt $values[$r].Count; $c++) {\n        $t.Cell($r + 1, $c + 1).Range.Text = $values[$r][$c]\n    }\n}\n"}
